$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    # Force Excel to store the literal text (not re-parse "5.0" -> 5) by using
    # the leading-apostrophe text-entry convention, then restore the cell's
    # original style so number formatting doesn't change.
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

function Set-Row($row, $a, $b, $c, $d, $e) {
    Set-TextCell $row 1 $a
    Set-TextCell $row 2 $b
    Set-TextCell $row 3 $c
    Set-TextCell $row 4 $d
    Set-TextCell $row 5 $e
}

# Update existing rows 2-6 with new values per the diff
Set-Row 2 "Natan" "5.0" "5.0" "5.0" "Em Recuperação"
Set-Row 3 "Natan" "7.0" "6.0" "6.5" "Em Recuperação"
Set-Row 4 "Teste" "10.0" "10.0" "10.0" "Aprovado"
Set-Row 5 "Tati" "5.0" "5.0" "5.0" "Em Recuperação"
Set-Row 6 "Natan" "10.0" "10.0" "10.0" "Aprovado"

# Add new rows 7-12
Set-Row 7 "Natan" "5.0" "5.0" "5.0" "Em Recuperação"
Set-Row 8 " " " " " " " " " "
Set-Row 9 "Natan" "5.0" "5.0" "5.0" "Em Recuperação"
Set-Row 10 " " " " " " " " " "
Set-Row 11 "Teste" "9.0" "9.0" "9.0" "Aprovado"
Set-Row 12 " " " " " " " " " "
